$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Agosto de 2020 a las 23:05"

# --- Update country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) for refreshed rows ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5464379
$ws.Range("C4").Value = 48713
$ws.Range("D4").Value = 2864431
$ws.Range("E4").Value = 2428665
$ws.Range("G4").Value = 868
$ws.Range("H4").Value = 171283

# Sudafrica (row 8)
$ws.Range("B8").Value = 579140
$ws.Range("C8").Value = 6275
$ws.Range("D8").Value = 461734
$ws.Range("E8").Value = 105850
$ws.Range("G8").Value = 286
$ws.Range("H8").Value = 11556

# Alemania (row 22)
$ws.Range("B22").Value = 223774
$ws.Range("C22").Value = 1505
$ws.Range("E22").Value = 13685
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = 9289

# Francia (row 23)
$ws.Range("D23").Value = 83848
$ws.Range("E23").Value = 97957

# Costa Rica (row 69)
$ws.Range("B69").Value = 26931
$ws.Range("C69").Value = 802
$ws.Range("D69").Value = 8785
$ws.Range("E69").Value = 17865
$ws.Range("G69").Value = 9
$ws.Range("H69").Value = 281

# Costa de Marfil (row 76)
$ws.Range("B76").Value = 16935
$ws.Range("C76").Value = 46
$ws.Range("D76").Value = 13721
$ws.Range("E76").Value = 3106
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 108

# Guayana Francesa (row 90)
$ws.Range("B90").Value = 8549
$ws.Range("C90").Value = 78
$ws.Range("D90").Value = 7841
$ws.Range("E90").Value = 655
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = 53

# Guinea (row 92)
$ws.Range("B92").Value = 8260
$ws.Range("C92").Value = 62
$ws.Range("D92").Value = 7177
$ws.Range("E92").Value = 1033

# Republica de Yibuti (row 105)
$ws.Range("B105").Value = 5367
$ws.Range("C105").Value = 9
$ws.Range("D105").Value = 5181
$ws.Range("E105").Value = 127

# Zimbabue (row 106)
$ws.Range("B106").Value = 5072
$ws.Range("C106").Value = 82
$ws.Range("D106").Value = 1998
$ws.Range("E106").Value = 2946

# --- Reorder + refresh Chad / Vietnam / Aruba block (rows 157-159) ---
# New rank order (by Casos totales desc): Aruba, Republica del Chad, Vietnam
$ws.Range("A157").Value = "Aruba"
$ws.Range("B157").Value = 973
$ws.Range("C157").Value = 79
$ws.Range("D157").Value = 114
$ws.Range("E157").Value = 855
$ws.Range("H157").Value = 4

$ws.Range("A158").Value = "Republica del Chad"
$ws.Range("B158").Value = 951
$ws.Range("C158").Value = 2
$ws.Range("D158").Value = 862
$ws.Range("E158").Value = 13
$ws.Range("H158").Value = 76

$ws.Range("A159").Value = "Vietnam"
$ws.Range("B159").Value = 929
$ws.Range("C159").Value = 18
$ws.Range("D159").Value = 437
$ws.Range("E159").Value = 471
$ws.Range("H159").Value = 21

# Reunion (row 162)
$ws.Range("B162").Value = 776
$ws.Range("C162").Value = 22
$ws.Range("D162").Value = 657
$ws.Range("E162").Value = 114

# Monaco (row 187)
$ws.Range("B187").Value = 146
$ws.Range("C187").Value = 2
$ws.Range("E187").Value = 28

# --- Reorder Islas Malvinas / Montserrat block (rows 213-214) ---
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
